$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "[1, 0, 0, 1, 0, 0, 1]"
$ws.Range("E3").Value = "['Normal', 'ParamViolation', 'SoftwareFault']"

# Row 6
$ws.Range("D6").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E6").Value = "['Normal', 'ParamViolation']"

# Row 11
$ws.Range("D11").Value = "[1, 0, 1, 0, 1, 0, 0]"
$ws.Range("E11").Value = "['Normal', 'HardwareFault', 'RegulationViolation']"

# Row 15
$ws.Range("D15").Value = "[0, 0, 0, 1, 0, 0, 0]"
$ws.Range("E15").Value = "['ParamViolation']"

# Row 16
$ws.Range("D16").Value = "[1, 0, 0, 0, 1, 0, 0]"
$ws.Range("E16").Value = "['Normal', 'RegulationViolation']"

# Row 24
$ws.Range("D24").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E24").Value = "[]"

# Row 27
$ws.Range("D27").Value = "[0, 0, 1, 0, 0, 0, 1]"
$ws.Range("E27").Value = "['HardwareFault', 'SoftwareFault']"

# Row 28
$ws.Range("D28").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E28").Value = "['Normal', 'SoftwareFault']"

# Row 29
$ws.Range("D29").Value = "[0, 0, 0, 1, 0, 0, 1]"
$ws.Range("E29").Value = "['ParamViolation', 'SoftwareFault']"

# Row 54
$ws.Range("D54").Value = "[0, 0, 0, 0, 0, 1, 0]"
$ws.Range("E54").Value = "['CommunicationIssue']"

# Row 56
$ws.Range("D56").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E56").Value = "['HardwareFault']"

# Row 68
$ws.Range("D68").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E68").Value = "['Normal', 'ParamViolation']"

# Row 71
$ws.Range("D71").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E71").Value = "['Normal', 'ParamViolation']"

# Row 73
$ws.Range("D73").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'ParamViolation']"

# Row 74
$ws.Range("D74").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E74").Value = "['Normal', 'SoftwareFault']"

# Row 89
$ws.Range("D89").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E89").Value = "['Normal', 'ParamViolation']"

# Row 109
$ws.Range("D109").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E109").Value = "['Normal', 'SurroundingEnvironment']"

# Row 113
$ws.Range("D113").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E113").Value = "['Normal', 'HardwareFault']"
